$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new job posting row (row 38): Taoyuan Airport summer internship
$ws.Range("H38").Value = "https://drive.google.com/drive/folders/1O2YeyP1td0ciqUunYCTD48k7vSYGG-XH?dmr=1&ec=wgc-drive-globalnav-goto"
$ws.Range("A38").Value = "桃園機場暑期實習"
$ws.Range("D38").Value = "桃園"
$ws.Range("E38").Value = "大三以上"
$ws.Range("B38").Value = "系辦申請"
$ws.Range("C38").Value = 45785

# Leave the cursor where the author left off after entering the new row
$ws.Range("E38").Select()
